# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows keyed by event name
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 601   # 合肥·Holic动漫游戏展
$wsExpo.Range("F6").Value = 360   # 合肥·W·A第五人格同人only2.0
$wsExpo.Range("F7").Value = 1835  # 合肥·第九届环形宇宙动漫游戏嘉年华
$wsExpo.Range("F8").Value = 100   # 合肥·MAX特摄同人only2.0

# Sheet "全部类型" (all categories) - same events, different row offsets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 601    # 合肥·Holic动漫游戏展
$wsAll.Range("F6").Value = 360    # 合肥·W·A第五人格同人only2.0
$wsAll.Range("F11").Value = 1835  # 合肥·第九届环形宇宙动漫游戏嘉年华
$wsAll.Range("F12").Value = 100   # 合肥·MAX特摄同人only2.0

$wb.Save()
